# "Generate Report for Archive"
# The localization status for this record has moved out of "Ready for
# handoff" and into "In Translation" — update every place that status is
# recorded: the per-language report sheets (zh-cn / de-de, column
# "Status") and the rollup "Overview" sheet (the zh-cn / de-de columns).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    # Track which columns actually received the new (shorter) status text so
    # only those get re-fit below — other columns' widths are untouched.
    $touchedCols = @()
    foreach ($row in 1..$used.Rows.Count) {
        foreach ($col in 1..$used.Columns.Count) {
            $cell = $ws.Cells.Item($row, $col)
            # Cast to [string] (and put the literal on the left) so PowerShell's
            # type-coercing -eq never turns this into a boolean comparison for
            # the True/False status cells elsewhere in the sheet.
            $cellText = [string]$cell.Value2
            if ($oldStatus -eq $cellText) {
                $cell.Value = $newStatus
                $touchedCols += $col
            }
        }
    }
    # The status column(s) narrow now that "In Translation" is shorter than
    # "Ready for handoff" — re-fit just those columns to their new content,
    # matching what Excel does on save. Columns with unrelated content keep
    # their existing width.
    foreach ($col in ($touchedCols | Select-Object -Unique)) {
        $ws.Columns.Item($col).AutoFit()
    }
}
